$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 15.62872253016054
$ws.Range("C2").Value = 10.97206875514792
$ws.Range("E2").Value = 23.88969109535498
$ws.Range("F2").Value = 39.82592705794573
$ws.Range("G2").Value = 25.07726275764634
$ws.Range("H2").Value = 13.18586315583154
$ws.Range("I2").Value = 17.72277532159539
$ws.Range("J2").Value = 7.663857844816457
$ws.Range("B3").Value = 14.85001136982871
$ws.Range("C3").Value = 10.25671156437281
$ws.Range("E3").Value = 23.60183170246736
$ws.Range("F3").Value = 39.62542244807391
$ws.Range("G3").Value = 25.17379221869006
$ws.Range("H3").Value = 13.27580412732864
$ws.Range("I3").Value = 17.9259863234315
$ws.Range("J3").Value = 7.701038358373495
$ws.Range("B4").Value = 14.35105579784221
$ws.Range("C4").Value = 9.802534016157381
$ws.Range("E4").Value = 23.42715513978223
$ws.Range("F4").Value = 39.51598614065925
$ws.Range("G4").Value = 25.25273898077442
$ws.Range("H4").Value = 13.33534938594793
$ws.Range("I4").Value = 18.05761427594968
$ws.Range("J4").Value = 7.725255727264973
$ws.Range("B5").Value = 14.14271337664081
$ws.Range("C5").Value = 9.615363806136864
$ws.Range("E5").Value = 23.35656468224052
$ws.Range("F5").Value = 39.47485866891193
$ws.Range("G5").Value = 25.28978599932364
$ws.Range("H5").Value = 13.36069574935698
$ws.Range("I5").Value = 18.1129756864495
$ws.Range("J5").Value = 7.735473776916225
$ws.Range("B6").Value = 14.10782269513724
$ws.Range("C6").Value = 9.583920540413821
$ws.Range("E6").Value = 23.34488098172634
$ws.Range("F6").Value = 39.46823974544179
$ws.Range("G6").Value = 25.29622987414223
$ws.Range("H6").Value = 13.36496963886638
$ws.Range("I6").Value = 18.12227235698197
$ws.Range("J6").Value = 7.737191576554614
$ws.Range("B7").Value = 14.34826599877582
$ws.Range("C7").Value = 9.800034206367508
$ws.Range("E7").Value = 23.42620064346348
$ws.Range("F7").Value = 39.51541740141168
$ws.Range("G7").Value = 25.25321896891759
$ws.Range("H7").Value = 13.33568684531633
$ws.Range("I7").Value = 18.05835392964152
$ws.Range("J7").Value = 7.725392117046284
$ws.Range("B8").Value = 15.36466689852039
$ws.Range("C8").Value = 10.73113941484857
$ws.Range("E8").Value = 23.79005041176231
$ws.Range("F8").Value = 39.75397495066562
$ws.Range("G8").Value = 25.1064203436458
$ws.Range("H8").Value = 13.21597467713905
$ws.Range("I8").Value = 17.79141836377424
$ws.Range("J8").Value = 7.676389617447153
$ws.Range("B9").Value = 17.18460746713687
$ws.Range("C9").Value = 12.36371520886841
$ws.Range("E9").Value = 24.51682878744608
$ws.Range("F9").Value = 40.32866664121408
$ws.Range("G9").Value = 24.97759491964827
$ws.Range("H9").Value = 13.01575120956964
$ws.Range("I9").Value = 17.32244981053342
$ws.Range("J9").Value = 7.59130512821139
$ws.Range("B10").Value = 18.40752784326486
$ws.Range("C10").Value = 13.43119624613443
$ws.Range("E10").Value = 25.05464145929913
$ws.Range("F10").Value = 40.81348072246748
$ws.Range("G10").Value = 24.98336997222192
$ws.Range("H10").Value = 12.8900408395083
$ws.Range("I10").Value = 17.01126391411006
$ws.Range("J10").Value = 7.535495671474464
$ws.Range("B11").Value = 18.93775795488071
$ws.Range("C11").Value = 13.88843278109591
$ws.Range("E11").Value = 25.29924879046634
$ws.Range("F11").Value = 41.04696450701441
$ws.Range("G11").Value = 25.0083762489916
$ws.Range("H11").Value = 12.83757338746085
$ws.Range("I11").Value = 16.87698195274308
$ws.Range("J11").Value = 7.511560258710674
$ws.Range("B12").Value = 19.13470827436996
$ws.Range("C12").Value = 14.05751860134321
$ws.Range("E12").Value = 25.39179507691565
$ws.Range("F12").Value = 41.13717648037725
$ws.Range("G12").Value = 25.02109775610104
$ws.Range("H12").Value = 12.81839042604793
$ws.Range("I12").Value = 16.82718376689859
$ws.Range("J12").Value = 7.502705395558404
$ws.Range("B13").Value = 19.09246331103429
$ws.Range("C13").Value = 14.02128319407013
$ws.Range("E13").Value = 25.37186835940546
$ws.Range("F13").Value = 41.11766887077396
$ws.Range("G13").Value = 25.01821286764853
$ws.Range("H13").Value = 12.8224912310275
$ws.Range("I13").Value = 16.83786184890942
$ws.Range("J13").Value = 7.504603153775339
$ws.Range("B14").Value = 18.95403849412775
$ws.Range("C14").Value = 13.90242482555253
$ws.Range("E14").Value = 25.30686467780005
$ws.Range("F14").Value = 41.05435067253769
$ws.Range("G14").Value = 25.00935751470596
$ws.Range("H14").Value = 12.83598141959448
$ws.Range("I14").Value = 16.87286392294353
$ws.Range("J14").Value = 7.510827575416883
$ws.Range("B15").Value = 18.86874743416448
$ws.Range("C15").Value = 13.8290925199544
$ws.Range("E15").Value = 25.26703528595321
$ws.Range("F15").Value = 41.01579846470557
$ws.Range("G15").Value = 25.00435773622098
$ws.Range("H15").Value = 12.84433400183832
$ws.Range("I15").Value = 16.89444079163945
$ws.Range("J15").Value = 7.514667428684361
$ws.Range("B16").Value = 18.3723435027261
$ws.Range("C16").Value = 13.4007456103532
$ws.Range("E16").Value = 25.03864845574547
$ws.Range("F16").Value = 40.7984770155357
$ws.Range("G16").Value = 24.98218920422526
$ws.Range("H16").Value = 12.89356524096795
$ws.Range("I16").Value = 17.02018641688262
$ws.Range("J16").Value = 7.537089142823626
$ws.Range("B17").Value = 18.06106901266788
$ws.Range("C17").Value = 13.13071710528637
$ws.Range("E17").Value = 24.89847515627588
$ws.Range("F17").Value = 40.66842734474772
$ws.Range("G17").Value = 24.97434854899977
$ws.Range("H17").Value = 12.92498038798724
$ws.Range("I17").Value = 17.09919470117226
$ws.Range("J17").Value = 7.55121622789921
$ws.Range("B18").Value = 17.8795811012832
$ws.Range("C18").Value = 12.972732749594
$ws.Range("E18").Value = 24.81785119381329
$ws.Range("F18").Value = 40.5948482545371
$ws.Range("G18").Value = 24.9719437835512
$ws.Range("H18").Value = 12.9434932956913
$ws.Range("I18").Value = 17.14532289909831
$ws.Range("J18").Value = 7.559478471885164
$ws.Range("B19").Value = 17.81771433036539
$ws.Range("C19").Value = 12.91878261319079
$ws.Range("E19").Value = 24.79055560329298
$ws.Range("F19").Value = 40.5701474020355
$ws.Range("G19").Value = 24.97148991208572
$ws.Range("H19").Value = 12.94983745787398
$ws.Range("I19").Value = 17.16105857479113
$ws.Range("J19").Value = 7.562299405285511
$ws.Range("B20").Value = 18.09445906231415
$ws.Range("C20").Value = 13.15973841026464
$ws.Range("E20").Value = 24.91339738690441
$ws.Range("F20").Value = 40.68214532744297
$ws.Range("G20").Value = 24.97496506759913
$ws.Range("H20").Value = 12.92159021750803
$ws.Range("I20").Value = 17.09071323462723
$ws.Range("J20").Value = 7.549698225311439
$ws.Range("B21").Value = 18.99480193629631
$ws.Range("C21").Value = 13.93744643531196
$ws.Range("E21").Value = 25.32596062467507
$ws.Range("F21").Value = 41.07290053007516
$ws.Range("G21").Value = 25.01187005825841
$ws.Range("H21").Value = 12.83200037019549
$ws.Range("I21").Value = 16.86255438691141
$ws.Range("J21").Value = 7.508993640080643
$ws.Range("B22").Value = 19.56083839169995
$ws.Range("C22").Value = 14.42206938895016
$ws.Range("E22").Value = 25.59508787123381
$ws.Range("F22").Value = 41.33872182182154
$ws.Range("G22").Value = 25.05495489382851
$ws.Range("H22").Value = 12.77744604356255
$ws.Range("I22").Value = 16.71957085541739
$ws.Range("J22").Value = 7.483608864072858
$ws.Range("B23").Value = 19.2608069130902
$ws.Range("C23").Value = 14.16557496098377
$ws.Range("E23").Value = 25.45152015172026
$ws.Range("F23").Value = 41.19591491592674
$ws.Range("G23").Value = 25.0302155539773
$ws.Range("H23").Value = 12.80619464768323
$ws.Range("I23").Value = 16.79532107474123
$ws.Range("J23").Value = 7.497045714910502
$ws.Range("B24").Value = 18.07937130101378
$ws.Range("C24").Value = 13.14662641858979
$ws.Range("E24").Value = 24.90665115434684
$ws.Range("F24").Value = 40.67593972152665
$ws.Range("G24").Value = 24.97467979184872
$ws.Range("H24").Value = 12.92312150637319
$ws.Range("I24").Value = 17.09454550907922
$ws.Range("J24").Value = 7.55038407697084
$ws.Range("B25").Value = 16.71180009634475
$ws.Range("C25").Value = 11.94532545952756
$ws.Range("E25").Value = 24.31923391612025
$ws.Range("F25").Value = 40.1619980372956
$ws.Range("G25").Value = 24.99503138413772
$ws.Range("H25").Value = 13.06618522581426
$ws.Range("I25").Value = 17.44347100764977
$ws.Range("J25").Value = 7.61314518528953
